# Apply crypto price/volume updates per the commit diff.
# Column D stores prices as text. A handful of the new values look
# numeric (e.g. "527.60", "4.00", "0.0000137") and Excel would silently
# reinterpret a plain Value assignment as a number, dropping the
# significant trailing/leading zeros. Prefix those with a leading
# apostrophe to force text entry, then reset .Style so the cell keeps
# its original (unstyled) formatting -- only the stored value changes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.874.19'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '2.502.82'
$ws.Range('E3').Value = '  +1.36%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = "'527.60"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.69%  '
$ws.Range('D6').Value = "'134.70"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.86%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('E8').Value = '  +1.61%  '
$ws.Range('D9').Value = '2.545.89'
$ws.Range('E9').Value = '  +2.67%  '
$ws.Range('D10').Value = "'0.0994"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.25%  '
$ws.Range('E11').Value = '  -1.36%  '
$ws.Range('E12').Value = '  -1.77%  '
$ws.Range('D13').Value = "'0.337"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.18%  '
$ws.Range('D14').Value = '2.966.82'
$ws.Range('E14').Value = '  +1.96%  '
$ws.Range('D15').Value = '58.770.59'
$ws.Range('E15').Value = '  +1.36%  '
$ws.Range('D16').Value = "'22.49"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.40%  '
$ws.Range('D17').Value = "'0.0000137"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.62%  '
$ws.Range('D18').Value = '2.536.05'
$ws.Range('E18').Value = '  +2.25%  '
$ws.Range('D19').Value = "'10.77"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.66%  '
$ws.Range('D20').Value = "'325.27"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.83%  '
$ws.Range('D21').Value = "'4.22"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.53%  '
$ws.Range('E22').Value = '  +7.11%  '
$ws.Range('D23').Value = "'0.995"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.34%  '
$ws.Range('D24').Value = "'65.34"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.16%  '
$ws.Range('D25').Value = "'0.412"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.59%  '
$ws.Range('E26').Value = '  +0.02%  '
$ws.Range('D27').Value = "'0.995"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.28%  '
$ws.Range('D28').Value = "'7.48"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.87%  '
$ws.Range('D29').Value = '0.0₃0762'
$ws.Range('E29').Value = '  +2.16%  '
$ws.Range('E30').Value = '  +3.08%  '
$ws.Range('E31').Value = '  +3.07%  '
$ws.Range('D32').Value = "'6.41"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.29%  '
$ws.Range('D33').Value = "'168.48"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.48%  '
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('D35').Value = "'0.991"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.63%  '
$ws.Range('D36').Value = "'18.40"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.95%  '
$ws.Range('E37').Value = '  -1.50%  '
$ws.Range('D38').Value = "'4.00"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.22%  '
$ws.Range('E39').Value = '  +3.21%  '
$ws.Range('D40').Value = "'36.74"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.36%  '
$ws.Range('D41').Value = "'0.793"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.37%  '
$ws.Range('D42').Value = "'284.45"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.81%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = "'5.19"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.35%  '
$ws.Range('B44').Value = 'Filecoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D44').Value = "'3.50"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').Value = "'130.97"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.18%  '
$ws.Range('D46').Value = "'0.607"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.22%  '
$ws.Range('D47').Value = "'0.0921"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.39%  '
$ws.Range('D48').Value = "'0.0509"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.15%  '
$ws.Range('D49').Value = "'17.94"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.84%  '
$ws.Range('E50').Value = '  +1.95%  '
$ws.Range('D51').Value = "'17.33"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.06%  '
